# Generate Report for Handoff
#
# Updates the localization-status report to reflect that a fresh handoff
# xliff generation pass has completed for the six files that were
# "Ready for handoff" (rows 7, 8, 9, 12, 13, 14 on the per-language
# sheets / the matching rows on the Overview sheet):
#   - Priority is now flagged "ht" on the zh-cn and de-de sheets
#   - The Latest Handoff Datetime is refreshed on zh-cn and de-de
#   - The Latest HO Xliff Generate Date is refreshed on the Overview sheet

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 12, 13, 14)

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-22 10:01:36"
}

# zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-22 10:01:31"
}

# de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-22 10:01:36"
}
